# Weekly refresh: insert a new latest-week pair of rows (Primera/Segunda)
# at the top of the data block for "Betarraga" at Terminal La Palmera de
# La Serena, pushing all the existing weekly rows down by two rows.

$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Insert two new blank rows before row 130 (the existing data there, and
# everything below it, shifts down to rows 132.. onward).
$ws.Rows.Item(130).Resize(2).Insert()

# Row 130: new "Primera" quality week record (2022-01-06).
$ws.Range("A130").Value = 8
$ws.Range("B130").Value = "Terminal La Palmera de La Serena"
$ws.Range("C130").Value = "Coquimbo"
$ws.Range("D130").Value = 44567
$ws.Range("E130").Value = 4
$ws.Range("F130").Value = 100114014
$ws.Range("G130").Value = "Betarraga"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 3060
$ws.Range("K130").Value = 450
$ws.Range("L130").Value = 500
$ws.Range("M130").Value = 475
$ws.Range("N130").Value = '$/paquete 3 unidades'
$ws.Range("O130").Value = "Provincia del Elquí"
$ws.Range("P130").Value = 158
$ws.Range("Q130").Value = 3
$ws.Range("R130").Value = "Hortaliza"

# Row 131: new "Segunda" quality week record (2022-01-06).
$ws.Range("A131").Value = 8
$ws.Range("B131").Value = "Terminal La Palmera de La Serena"
$ws.Range("C131").Value = "Coquimbo"
$ws.Range("D131").Value = 44567
$ws.Range("E131").Value = 4
$ws.Range("F131").Value = 100114014
$ws.Range("G131").Value = "Betarraga"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Segunda"
$ws.Range("J131").Value = 1520
$ws.Range("K131").Value = 350
$ws.Range("L131").Value = 400
$ws.Range("M131").Value = 375
$ws.Range("N131").Value = '$/paquete 3 unidades'
$ws.Range("O131").Value = "Provincia del Elquí"
$ws.Range("P131").Value = 125
$ws.Range("Q131").Value = 3
$ws.Range("R131").Value = "Hortaliza"
